$d = $word.ActiveDocument

# -----------------------------------------------------------------
# 1) Merge the three runs describing the k=2 cluster stability test
#    back into a single run (no visible text change, Find/Replace
#    with identical text collapses adjacent equal-format runs).
# -----------------------------------------------------------------
$bodyRng = $d.Content
$bodyRng.Find.Execute(
    "The k=2 clusters were also found to be highly stable; when the analysis was repeated on 5 sets of perturbed observations, for all sets ",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "The k=2 clusters were also found to be highly stable; when the analysis was repeated on 5 sets of perturbed observations, for all sets ",
    2) | Out-Null

# -----------------------------------------------------------------
# 2) Insert two new paragraphs (one blank, one with new body text)
#    plus a trailing blank paragraph after the SLC45A2 sentence,
#    just before the blank "keepNext" paragraph leading into
#    REFERENCES.
# -----------------------------------------------------------------
$findRng = $d.Content
$findRng.Find.Execute("EAS super-population.", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$findRng.Collapse(0)
$findRng.InsertParagraphAfter()
$findRng.Collapse(0)
$findRng.InsertParagraphAfter()
$findRng.Collapse(0)
$findRng.InsertParagraphAfter()

# Re-locate the three freshly inserted (currently empty) paragraphs.
$anchorRng = $d.Content
$anchorRng.Find.Execute("EAS super-population.", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$anchorPara = $anchorRng.Paragraphs(1)
$blankPara1 = $anchorPara.Next()
$textPara   = $blankPara1.Next()
$blankPara2 = $textPara.Next()

$textRng = $textPara.Range
$textRng.ParagraphFormat.FirstLineIndent = 36
$textRng.InsertBefore("Decision trees with bagging are a very easy and effective way to perform classification as they use a simple yes/no format to classify new data into given groups. This is a very intuitive way to classify, especially genotype data. Genotypes are one of three possibilities: homozygous dominant, heterozygous, or homozygous recessive. This provides the decision tree with very clear boundaries for branches and makes identification of influential SNPs more straight forward. Unfortunately, with very large datasets, which is common with genotypes, bagging is very computationally expensive. This being said, bagging provided good results and in theory would be a great application for population clustering based on SNPs. Overall, decision trees and bagging are a great tool, but logistic regression and other clustering techniques are able to identify similar patterns within the data, with far less computational effort.  ")
$textRng.Font.NameBi = "Calibri"

# Split "easy" into its own run (matching the original authoring history)
# by toggling a character attribute off/on; this forces the run boundary
# without altering the final visible formatting.
$wordSearchRng = $d.Range($textPara.Range.Start, $textPara.Range.End)
$wordSearchRng.Find.Execute("easy", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$wordSearchRng.Bold = 1
$wordSearchRng.Bold = 0

# -----------------------------------------------------------------
# 3) Merge the author-name runs in the header (", Heather ",
#    "Treleaven", ", Natalie Kim") into a single run, dropping the
#    now-unneeded spell-check markers around "Treleaven".
# -----------------------------------------------------------------
$headers = $d.Sections(1).Headers
$primaryHeader = $headers.Item(2)
$headerRng = $primaryHeader.Range
$headerRng.Find.Execute(
    ", Heather Treleaven, Natalie Kim",
    $true, $false, $false, $false, $false, $true, 1, $false,
    ", Heather Treleaven, Natalie Kim",
    2) | Out-Null
